$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.377.29"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.34"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.89"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4493"
$ws.Range("E7").Value = "  +5.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3783"
$ws.Range("E8").Value = "  +3.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07512"
$ws.Range("E9").Value = "  +4.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8976"
$ws.Range("E10").Value = "  +6.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.10"
$ws.Range("E11").Value = "  +1.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.817.48"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.784"
$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.64"
$ws.Range("E14").Value = "  +5.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.414"
$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9987"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.23"
$ws.Range("E20").Value = "  +2.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.389.48"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.302"
$ws.Range("E22").Value = "  +3.38%  "

$ws.Range("E23").Value = "  +1.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.051.30"
$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.005"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("E26").Value = "  +10.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.65"
$ws.Range("E28").Value = "  +2.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.400"
$ws.Range("E29").Value = "  +3.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.28"
$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08847"
$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7812"
$ws.Range("E32").Value = "  +6.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.201"
$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.594"
$ws.Range("E34").Value = "  +4.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.883"
$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9989"
$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  +2.07%  "

$ws.Range("E38").Value = "  +2.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05341"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.434"
$ws.Range("E40").Value = "  +2.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5367"
$ws.Range("E41").Value = "  +4.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1736"
$ws.Range("E42").Value = "  +2.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.861"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.281"
$ws.Range("E44").Value = "  +16.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.845"
$ws.Range("E45").Value = "  +3.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5166"
$ws.Range("E46").Value = "  +9.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.83"
$ws.Range("E47").Value = "  +3.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.93"
$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.710"
$ws.Range("E49").Value = "  +3.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9986"
$ws.Range("E50").Value = "  -0.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06388"
$ws.Range("E51").Value = "  +1.02%  "

